$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 15 more test cities (Redfin test data) to column A, rows 49-63.
# NOTE: the shared-strings table must be populated in the exact order the
# original author typed/pasted them, which is not strictly row order
# ("Los Angeles", destined for row 51, was entered/registered before
# "San jose", destined for row 50). We therefore set the cell values in
# that same creation order while still targeting the correct row/cell
# for each city, so both the sheet data and the shared string table
# line up with the target workbook.

$ws.Cells.Item(49, 1).Value = "Virginia Beach"
$ws.Cells.Item(51, 1).Value = "Los Angeles"
$ws.Cells.Item(50, 1).Value = "San jose"
$ws.Cells.Item(52, 1).Value = "Philadelphia"
$ws.Cells.Item(53, 1).Value = "Indianapolis"
$ws.Cells.Item(54, 1).Value = "El Paso"
$ws.Cells.Item(55, 1).Value = "Tucson"
$ws.Cells.Item(56, 1).Value = "Omaha"
$ws.Cells.Item(57, 1).Value = "Wichita"
$ws.Cells.Item(58, 1).Value = "Cleveland"
$ws.Cells.Item(59, 1).Value = "Cincinnati"
$ws.Cells.Item(60, 1).Value = "Pittsburgh"
$ws.Cells.Item(61, 1).Value = "Toldedo"
$ws.Cells.Item(62, 1).Value = "Norfolk"
$ws.Cells.Item(63, 1).Value = "Salt Lake City"

# Mirror the author's final selection/scroll state: the active cell moves
# to the next blank row below the newly-added data.
$ws.Range("A64").Select() | Out-Null
